# Auto-generated edit script: updates cached market-price derived values
# across the Leve profit tables (Table_ALC .. Table_WVR) to match the
# refreshed scheduled-runner data pull described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 794418.0600000001
$ws.Range("I28").Value = 1234716.4
$ws.Range("J28").Value = 1881.2
$ws.Range("K28").Value = 1234716.4
$ws.Range("L28").Value = 1881.2
$ws.Range("M28").Value = -1234231.4
$ws.Range("N28").Value = -2851.2
$ws.Range("H40").Value = 2330
$ws.Range("J40").Value = 2220
$ws.Range("L40").Value = 2220
$ws.Range("N40").Value = -2570
$ws.Range("H55").Value = 233.63637
$ws.Range("I55").Value = 248
$ws.Range("J55").Value = 90
$ws.Range("K55").Value = 248
$ws.Range("L55").Value = 90
$ws.Range("M55").Value = -34
$ws.Range("N55").Value = -518
$ws.Range("H61").Value = 933.3333
$ws.Range("I61").Value = 933.3333
$ws.Range("K61").Value = 2799.9999
$ws.Range("M61").Value = -2627.9999
$ws.Range("H62").Value = 6958282.5
$ws.Range("I62").Value = 11121152
$ws.Range("K62").Value = 11121152
$ws.Range("M62").Value = -11120528
$ws.Range("H65").Value = 6958282.5
$ws.Range("I65").Value = 11121152
$ws.Range("K65").Value = 55605760
$ws.Range("M65").Value = -55602640
$ws.Range("H107").Value = 529634.75
$ws.Range("I107").Value = 794234.3
$ws.Range("J107").Value = 435.7143
$ws.Range("K107").Value = 794234.3
$ws.Range("L107").Value = 435.7143
$ws.Range("M107").Value = -792314.3
$ws.Range("N107").Value = -4275.7143
$ws.Range("H113").Value = 3000
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -9508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 47670.363
$ws.Range("I2").Value = 74148.07000000001
$ws.Range("K2").Value = 74148.07000000001
$ws.Range("M2").Value = -74035.07000000001
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H70").Value = 50577
$ws.Range("J70").Value = 50577
$ws.Range("L70").Value = 50577
$ws.Range("N70").Value = -51117
$ws.Range("H73").Value = 50577
$ws.Range("J73").Value = 50577
$ws.Range("L73").Value = 50577
$ws.Range("N73").Value = -52449
$ws.Range("H104").Value = 28633.334
$ws.Range("J104").Value = 28633.334
$ws.Range("L104").Value = 28633.334
$ws.Range("N104").Value = -35621.334
$ws.Range("H116").Value = 47670.363
$ws.Range("I116").Value = 74148.07000000001
$ws.Range("K116").Value = 74148.07000000001
$ws.Range("M116").Value = -71854.07000000001
$ws.Range("H132").Value = 2620.5789
$ws.Range("J132").Value = 4170.5454
$ws.Range("L132").Value = 12511.6362
$ws.Range("N132").Value = -17571.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 47670.363
$ws.Range("I3").Value = 74148.07000000001
$ws.Range("K3").Value = 74148.07000000001
$ws.Range("M3").Value = -74034.07000000001
$ws.Range("H5").Value = 13358.75
$ws.Range("I5").Value = 13358.75
$ws.Range("K5").Value = 13358.75
$ws.Range("M5").Value = -13245.75
$ws.Range("H107").Value = 1637.4286
$ws.Range("I107").Value = 1344.2222
$ws.Range("K107").Value = 1344.2222
$ws.Range("M107").Value = 575.7778000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1102
$ws.Range("I2").Value = 1102
$ws.Range("K2").Value = 1102
$ws.Range("M2").Value = -989
$ws.Range("H31").Value = 15382.556
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15382.556
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15382.556
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -15972.556
$ws.Range("H34").Value = 15382.556
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15382.556
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15382.556
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -15786.556
$ws.Range("H56").Value = 14000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3875
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5685
$ws.Range("H73").Value = 3875
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -4908
$ws.Range("H97").Value = 595.2857
$ws.Range("I97").Value = 695.75
$ws.Range("K97").Value = 2087.25
$ws.Range("M97").Value = -1591.25
$ws.Range("H102").Value = 4846.5713
$ws.Range("I102").Value = 3126
$ws.Range("J102").Value = 5133.3335
$ws.Range("K102").Value = 9378
$ws.Range("L102").Value = 15400.0005
$ws.Range("M102").Value = -6944
$ws.Range("N102").Value = -20268.0005
$ws.Range("H114").Value = 1027.1538
$ws.Range("I114").Value = 1275.5
$ws.Range("J114").Value = 916.7778
$ws.Range("K114").Value = 3826.5
$ws.Range("L114").Value = 2750.3334
$ws.Range("M114").Value = -572.5
$ws.Range("N114").Value = -9258.3334
$ws.Range("H115").Value = 1460
$ws.Range("I115").Value = 220
$ws.Range("J115").Value = 1666.6666
$ws.Range("K115").Value = 660
$ws.Range("L115").Value = 4999.9998
$ws.Range("M115").Value = 515
$ws.Range("N115").Value = -7349.9998
$ws.Range("H122").Value = 871.2857
$ws.Range("J122").Value = 1188.6666
$ws.Range("L122").Value = 10697.9994
$ws.Range("N122").Value = -15597.9994
$ws.Range("H132").Value = 10417633
$ws.Range("J132").Value = 15152657
$ws.Range("L132").Value = 136373913
$ws.Range("N132").Value = -136378973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66.333336
$ws.Range("I2").Value = 110
$ws.Range("J2").Value = 57.6
$ws.Range("K2").Value = 110
$ws.Range("L2").Value = 57.6
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = -283.6
$ws.Range("H4").Value = 28750
$ws.Range("J4").Value = 28750
$ws.Range("L4").Value = 28750
$ws.Range("N4").Value = -28974
$ws.Range("H80").Value = 2487.8096
$ws.Range("J80").Value = 3212.5
$ws.Range("L80").Value = 3212.5
$ws.Range("N80").Value = -5208.5
$ws.Range("H83").Value = 2487.8096
$ws.Range("J83").Value = 3212.5
$ws.Range("L83").Value = 16062.5
$ws.Range("N83").Value = -26046.5
$ws.Range("H113").Value = 2560
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9332.666999999999
$ws.Range("I22").Value = 925
$ws.Range("J22").Value = 13536.5
$ws.Range("K22").Value = 925
$ws.Range("L22").Value = 13536.5
$ws.Range("M22").Value = -630
$ws.Range("N22").Value = -14126.5
$ws.Range("H27").Value = 9332.666999999999
$ws.Range("I27").Value = 925
$ws.Range("J27").Value = 13536.5
$ws.Range("K27").Value = 925
$ws.Range("L27").Value = 13536.5
$ws.Range("M27").Value = -818
$ws.Range("N27").Value = -13750.5
$ws.Range("H40").Value = 3250.24
$ws.Range("I40").Value = 2875.6
$ws.Range("K40").Value = 2875.6
$ws.Range("M40").Value = -2739.6
$ws.Range("H61").Value = 1515.05
$ws.Range("I61").Value = 1162.5625
$ws.Range("K61").Value = 1162.5625
$ws.Range("M61").Value = -960.5625
$ws.Range("H68").Value = 2383.6667
$ws.Range("I68").Value = 1834
$ws.Range("J68").Value = 2933.3333
$ws.Range("K68").Value = 1834
$ws.Range("L68").Value = 2933.3333
$ws.Range("M68").Value = -1085
$ws.Range("N68").Value = -4431.3333
$ws.Range("H71").Value = 2383.6667
$ws.Range("I71").Value = 1834
$ws.Range("J71").Value = 2933.3333
$ws.Range("K71").Value = 9170
$ws.Range("L71").Value = 14666.6665
$ws.Range("M71").Value = -5426
$ws.Range("N71").Value = -22154.6665
$ws.Range("H106").Value = 23053.223
$ws.Range("J106").Value = 23053.223
$ws.Range("L106").Value = 23053.223
$ws.Range("N106").Value = -25577.223
$ws.Range("H113").Value = 1515.05
$ws.Range("I113").Value = 1162.5625
$ws.Range("K113").Value = 1162.5625
$ws.Range("M113").Value = 1007.4375
$ws.Range("H132").Value = 3905.6365
$ws.Range("I132").Value = 2805.6428
$ws.Range("K132").Value = 8416.928400000001
$ws.Range("M132").Value = -5886.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 12966.889
$ws.Range("J74").Value = 12590.429
$ws.Range("L74").Value = 12590.429
$ws.Range("N74").Value = -14462.429
$ws.Range("H77").Value = 12966.889
$ws.Range("J77").Value = 12590.429
$ws.Range("L77").Value = 37771.287
$ws.Range("N77").Value = -47131.287
$ws.Range("H107").Value = 644.7778
$ws.Range("I107").Value = 600.625
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 1801.875
$ws.Range("L107").Value = 2994
$ws.Range("M107").Value = 118.125
$ws.Range("N107").Value = -6834
$ws.Range("H112").Value = 30096.75
$ws.Range("J112").Value = 30096.75
$ws.Range("L112").Value = 30096.75
$ws.Range("N112").Value = -33050.75
$ws.Range("H113").Value = 413.63635
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 508.33334
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 1525.00002
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5865.000019999999
$ws.Range("H126").Value = 91407.82000000001
$ws.Range("I126").Value = 91407.82000000001
$ws.Range("K126").Value = 274223.46
$ws.Range("M126").Value = -271753.46
